$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells retain their original text (string) representation instead of
# being auto-converted to numbers by Excel, since the source data stores prices as text.
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D24","D25","D27","D40","D41","D42","D43","D44","D45","D47","D48","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values scraped on 2022-12-17 16:44 UTC
$ws.Range('D2').Value = '237.43'
$ws.Range('D3').Value = '22.04'
$ws.Range('D5').Value = '0.05596'
$ws.Range('D6').Value = '6.476'
$ws.Range('D7').Value = '3.338'
$ws.Range('D9').Value = '1.043'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').Value = '0.01167'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1384'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07300'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').Value = '0.03152'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = '0.02971'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '0.09237'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '0.001677'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'MCDex'
$ws.Range('C17').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D17').Value = '3.254'
$ws.Range('E17').Value = '16MCDexMCB'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').Value = '0.04769'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('D19').Value = '0.006219'
$ws.Range('D20').Value = '0.005059'
$ws.Range('D21').Value = '0.001049'
$ws.Range('D22').Value = '0.0001500'
$ws.Range('D24').Value = '3.972'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('D25').Value = '2.205'
$ws.Range('D27').Value = '0.1062'
$ws.Range('D40').Value = '0.04100'
$ws.Range('D41').Value = '0.007012'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '0.003502'
$ws.Range('E42').Value = '41CEJICEJI'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').Value = '0.1036'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').Value = '0.008812'
$ws.Range('D45').Value = '0.00005431'
$ws.Range('D47').Value = '0.6755'
$ws.Range('D48').Value = '0.03647'
$ws.Range('E48').Value = '47BOLOBOLOWorstin24h'
$ws.Range('D50').Value = '0.01010'
